$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$paraIndexes = @(1, 3, 5)
foreach ($i in $paraIndexes) {
    $para = $tr.Paragraphs($i)
    $para.Font.Size = 22
}
